$d = $word.ActiveDocument

# Locate the paragraph that starts the "Ver que hacer con la función ..." item
# (the skill_level / classif_origin warning) so the deletion isn't tied to a
# brittle hard-coded paragraph index.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Ver que hacer con la función*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Remove this paragraph and everything after it (the trailing blank
    # list paragraph included) through the end of the document body.
    $start = $target.Range.Start
    $end = $d.Content.End
    $r = $d.Range($start, $end)
    $r.Delete()
}
